# Update "Bestand" (column D) quantities on the "Lagerbestand M0129" sheet.
# These figures were revised during the config-file / common-function split
# described in the commit message; only the stock values change, all other
# cell content (labels, article numbers, styles) stays untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Lagerbestand M0129")

$updates = @{
    "D2"   = 475
    "D9"   = 2071
    "D14"  = 2248
    "D22"  = 1723
    "D29"  = 3161
    "D30"  = -991
    "D33"  = -97
    "D40"  = -62
    "D41"  = 2617
    "D42"  = 915
    "D43"  = -1952
    "D45"  = 2344
    "D46"  = 67
    "D48"  = -12
    "D51"  = -1240
    "D52"  = -8
    "D53"  = 28
    "D56"  = 1198
    "D63"  = -163
    "D75"  = -419
    "D76"  = 394
    "D82"  = 132
    "D90"  = 308
    "D91"  = 439
    "D92"  = 689
    "D93"  = 335
    "D97"  = 448
    "D99"  = 127
    "D118" = -166
    "D124" = 182
    "D125" = -64
    "D127" = 361
    "D155" = -1958
    "D160" = 408
    "D181" = -46
    "D182" = -124
    "D207" = -6229
    "D225" = -42
    "D230" = -1550
    "D231" = -22
    "D238" = -6200
    "D242" = -220
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
